$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 118, shifting existing rows 118:223 down to 119:224
$ws.Rows(118).Insert()

# Populate the newly inserted row 118 with the new record's data
$ws.Cells.Item(118, 1).Value = 11
$ws.Cells.Item(118, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(118, 3).Value = "Bíobío"
$ws.Cells.Item(118, 4).Value = 44705
$ws.Cells.Item(118, 5).Value = 8
$ws.Cells.Item(118, 6).Value = 100114013
$ws.Cells.Item(118, 7).Value = "Zanahoria"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 600
$ws.Cells.Item(118, 11).Value = 5500
$ws.Cells.Item(118, 12).Value = 6000
$ws.Cells.Item(118, 13).Value = 5750
$ws.Cells.Item(118, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(118, 15).Value = "Región de Ñuble"
$ws.Cells.Item(118, 16).Value = 288
$ws.Cells.Item(118, 17).Value = 20
$ws.Cells.Item(118, 18).Value = "Hortaliza"
